# "Generate Report for Handback" -- refresh the localization-status report:
#   * Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     on the Overview sheet (zh-cn/de-de columns) and on each language sheet's
#     Status column.
#   * The "Latest Handback DateTime" for each language's first file is bumped
#     to the new handback timestamp.
#   * The stale "handback file is not the latest" Error Detail is cleared now
#     that the handback is in sync.
#   * A couple of columns are widened/narrowed to fit the new text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.14437166849777
$wsOverview.Columns.Item(6).ColumnWidth = 29.14437166849777

# ---- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-10-13 14:11:44"
$wsZhCn.Range("K3").Value = "2016-10-13 14:11:44"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Columns.Item(3).ColumnWidth = 29.14437166849777
$wsZhCn.Columns.Item(16).ColumnWidth = 12.913719813028965

# ---- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-10-13 14:12:00"
$wsDeDe.Range("K3").Value = "2016-10-13 14:12:00"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = 29.14437166849777
$wsDeDe.Columns.Item(16).ColumnWidth = 12.913719813028965
